$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Shield",
    "Bounce Heal",
    "Unnamed AOE Bubble",
    "Regeneration",
    "Chloroplast",
    "Nature's Touch",
    "Replenishing Winds",
    "Blessing of the Grove Aura",
    "Blessing of the Grove HOT",
    "Blessing of the Grove Ref",
    "Replenish",
    "Living Seed",
    "Hibernate Friend",
    "Hibernate Foe"
)

$row = 16
$num = 15
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $num
    $ws.Cells.Item($row, 2).Value = $name
    $row++
    $num++
}

$ws.Columns.Item(2).ColumnWidth = 23.8

$ws.Range("A30").Select()
$excel.ActiveWindow.ScrollRow = 7
